$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 93837.37
$ws.Range("I70").Value = 1551
$ws.Range("J70").Value = 339934.34
$ws.Range("K70").Value = 4653
$ws.Range("L70").Value = 1019803.02
$ws.Range("M70").Value = -4383
$ws.Range("N70").Value = -1020343.02
$ws.Range("H73").Value = 93837.37
$ws.Range("I73").Value = 1551
$ws.Range("J73").Value = 339934.34
$ws.Range("K73").Value = 4653
$ws.Range("L73").Value = 1019803.02
$ws.Range("M73").Value = -3717
$ws.Range("N73").Value = -1021675.02
$ws.Range("H74").Value = 4550.375
$ws.Range("I74").Value = 4567.1665
$ws.Range("K74").Value = 4567.1665
$ws.Range("M74").Value = -3631.1665
$ws.Range("H77").Value = 4550.375
$ws.Range("I77").Value = 4567.1665
$ws.Range("K77").Value = 22835.8325
$ws.Range("M77").Value = -18155.8325
$ws.Range("H113").Value = 1908.3334
$ws.Range("I113").Value = 2140
$ws.Range("K113").Value = 2140
$ws.Range("M113").Value = 1114
$ws.Range("H135").Value = 19232538
$ws.Range("I135").Value = 705.25
$ws.Range("K135").Value = 6347.25
$ws.Range("M135").Value = -3812.25
$ws.Range("H137").Value = 374651.25
$ws.Range("I137").Value = 439913.12
$ws.Range("J137").Value = 4834
$ws.Range("K137").Value = 1319739.36
$ws.Range("L137").Value = 14502
$ws.Range("M137").Value = -1317189.36
$ws.Range("N137").Value = -19602

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4681.5835
$ws.Range("I32").Value = 5001.544
$ws.Range("K32").Value = 5001.544
$ws.Range("M32").Value = -4714.544
$ws.Range("H74").Value = 4143.154
$ws.Range("I74").Value = 885.4
$ws.Range("J74").Value = 5266.517
$ws.Range("K74").Value = 885.4
$ws.Range("L74").Value = 5266.517
$ws.Range("M74").Value = -11.39999999999998
$ws.Range("N74").Value = -7014.517
$ws.Range("H77").Value = 4143.154
$ws.Range("I77").Value = 885.4
$ws.Range("J77").Value = 5266.517
$ws.Range("K77").Value = 4427
$ws.Range("L77").Value = 26332.585
$ws.Range("M77").Value = -59
$ws.Range("N77").Value = -35068.585

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1432
$ws.Range("I64").Value = 3335.4285
$ws.Range("J64").Value = 730.7368
$ws.Range("K64").Value = 3335.4285
$ws.Range("L64").Value = 730.7368
$ws.Range("M64").Value = -3110.4285
$ws.Range("N64").Value = -1180.7368
$ws.Range("H67").Value = 1432
$ws.Range("I67").Value = 3335.4285
$ws.Range("J67").Value = 730.7368
$ws.Range("K67").Value = 3335.4285
$ws.Range("L67").Value = 730.7368
$ws.Range("M67").Value = -2555.4285
$ws.Range("N67").Value = -2290.7368
$ws.Range("H86").Value = 1975.2667
$ws.Range("I86").Value = 1677.4166
$ws.Range("J86").Value = 3166.6667
$ws.Range("K86").Value = 1677.4166
$ws.Range("L86").Value = 3166.6667
$ws.Range("M86").Value = -554.4166
$ws.Range("N86").Value = -5412.6667
$ws.Range("H89").Value = 1975.2667
$ws.Range("I89").Value = 1677.4166
$ws.Range("J89").Value = 3166.6667
$ws.Range("K89").Value = 8387.083000000001
$ws.Range("L89").Value = 15833.3335
$ws.Range("M89").Value = -2771.083000000001
$ws.Range("N89").Value = -27065.3335
$ws.Range("H105").Value = 2419.182
$ws.Range("I105").Value = 2957.1428
$ws.Range("J105").Value = 1477.75
$ws.Range("K105").Value = 2957.1428
$ws.Range("L105").Value = 1477.75
$ws.Range("M105").Value = -1210.1428
$ws.Range("N105").Value = -4971.75
$ws.Range("H107").Value = 232183.53
$ws.Range("I107").Value = 321313.12
$ws.Range("J107").Value = 1495.1177
$ws.Range("K107").Value = 321313.12
$ws.Range("L107").Value = 1495.1177
$ws.Range("M107").Value = -319393.12
$ws.Range("N107").Value = -5335.1177

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1661.6
$ws.Range("I31").Value = 1250.8334
$ws.Range("K31").Value = 1250.8334
$ws.Range("M31").Value = -955.8334
$ws.Range("H34").Value = 1661.6
$ws.Range("I34").Value = 1250.8334
$ws.Range("K34").Value = 1250.8334
$ws.Range("M34").Value = -1048.8334
$ws.Range("H58").Value = 3223.5642
$ws.Range("I58").Value = 1844.7858
$ws.Range("J58").Value = 3995.68
$ws.Range("K58").Value = 1844.7858
$ws.Range("L58").Value = 3995.68
$ws.Range("M58").Value = -1641.7858
$ws.Range("N58").Value = -4401.68
$ws.Range("H62").Value = 2651.6667
$ws.Range("I62").Value = 2582
$ws.Range("K62").Value = 2582
$ws.Range("M62").Value = -1958
$ws.Range("H65").Value = 2651.6667
$ws.Range("I65").Value = 2582
$ws.Range("K65").Value = 12910
$ws.Range("M65").Value = -9790
$ws.Range("H136").Value = 3223.5642
$ws.Range("I136").Value = 1844.7858
$ws.Range("J136").Value = 3995.68
$ws.Range("K136").Value = 5534.357400000001
$ws.Range("L136").Value = 11987.04
$ws.Range("M136").Value = -2984.357400000001
$ws.Range("N136").Value = -17087.04

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 41000
$ws.Range("J93").Value = 2000
$ws.Range("L93").Value = 6000
$ws.Range("N93").Value = -9744
$ws.Range("H97").Value = 2584
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 3040.8
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 9122.400000000001
$ws.Range("M97").Value = -404
$ws.Range("N97").Value = -10114.4
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").Value = ""
$ws.Range("H117").Value = 2803.4443
$ws.Range("I117").Value = 1040
$ws.Range("J117").Value = 5007.75
$ws.Range("K117").Value = 3120
$ws.Range("L117").Value = 15023.25
$ws.Range("M117").Value = 322
$ws.Range("N117").Value = -21907.25
$ws.Range("H129").Value = 1448
$ws.Range("I129").Value = 490
$ws.Range("J129").Value = 1687.5
$ws.Range("K129").Value = 1470
$ws.Range("L129").Value = 5062.5
$ws.Range("M129").Value = 3530
$ws.Range("N129").Value = -15062.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 464152.62
$ws.Range("J21").Value = 2001328
$ws.Range("L21").Value = 2001328
$ws.Range("N21").Value = -2001674
$ws.Range("H30").Value = 464152.62
$ws.Range("J30").Value = 2001328
$ws.Range("L30").Value = 2001328
$ws.Range("N30").Value = -2001538
$ws.Range("H70").Value = 5918.5
$ws.Range("I70").Value = 6286.857
$ws.Range("J70").Value = 5402.8
$ws.Range("K70").Value = 6286.857
$ws.Range("L70").Value = 5402.8
$ws.Range("M70").Value = -6016.857
$ws.Range("N70").Value = -5942.8
$ws.Range("H73").Value = 5918.5
$ws.Range("I73").Value = 6286.857
$ws.Range("J73").Value = 5402.8
$ws.Range("K73").Value = 6286.857
$ws.Range("L73").Value = 5402.8
$ws.Range("M73").Value = -5350.857
$ws.Range("N73").Value = -7274.8
$ws.Range("H122").Value = 5524.278
$ws.Range("I122").Value = 6680.0713
$ws.Range("J122").Value = 1479
$ws.Range("K122").Value = 20040.2139
$ws.Range("L122").Value = 4437
$ws.Range("M122").Value = -17590.2139
$ws.Range("N122").Value = -9337

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 104322.91
$ws.Range("I132").Value = 188009
$ws.Range("J132").Value = 3899.6
$ws.Range("K132").Value = 564027
$ws.Range("L132").Value = 11698.8
$ws.Range("M132").Value = -561497
$ws.Range("N132").Value = -16758.8

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1782.9546
$ws.Range("I132").Value = 1071.9286
$ws.Range("J132").Value = 3027.25
$ws.Range("K132").Value = 3215.7858
$ws.Range("L132").Value = 9081.75
$ws.Range("M132").Value = -685.7857999999997
$ws.Range("N132").Value = -14141.75
